$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts First Name/Last Name/Email/Phonenumber
# one column to the right (A->B, B->C, C->D, D->E) and carries their existing styling along.
$ws.Columns("A").Insert()

# Rename headers to match the new import format
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "first_name"
$ws.Range("C1").Value = "last_name"
$ws.Range("D1").Value = "email"
$ws.Range("E1").Value = "phone"

# Fill in the new "id" column
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4

# Update phone numbers to E.164-style formatted text (keep as text so the leading "+" is preserved)
$ws.Range("E2").Value = "+17068185081"
$ws.Range("E3").Value = "+14789733746"
$ws.Range("E4").Value = "+17063729685"
$ws.Range("E5").Value = "+19123348897"

# Restore selection to match the saved workbook state
$ws.Range("D9").Select()
